$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DOMA-3100: add a ":formatN()" carbone.io formatter to every ticket-status
# placeholder (processing/completed/canceled/deferred/closed/new_or_reopened)
# in the two data rows, so the exported values get converted to numbers.
# The address column (A) is left untouched.
$cols = @("B", "C", "D", "E", "F", "G")

foreach ($col in $cols) {
    foreach ($row in @(2, 3)) {
        $cell = $ws.Range($col + $row)
        $text = [string]$cell.Value2
        if ($text.EndsWith("}")) {
            $cell.Value2 = $text.Substring(0, $text.Length - 1) + ":formatN()}"
        }
    }
}

# Switch those same cells to a numeric (integer) format so the formatted
# numbers render correctly instead of as plain text.
$ws.Range("B2:G2").NumberFormat = "0"
$ws.Range("B3:G3").NumberFormat = "0"
